# "Generate Report for handback" - mark the xinjiang localization-status
# report as handed back: update the Status column, populate the
# Latest Target File / Latest Handback File columns with hyperlinks,
# stamp the Latest Handback DateTime, and flip the Handoff Reason to
# "Include" for the zh-cn and de-de hand-off rows.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# The "Status" column shows the same text on the Overview summary sheet
# (columns B/C) and on each language sheet (column B) - update them all.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$XlfFileName,
        [string]$XlfHyperlinkUrl,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status
    $ws.Range("B2").Value = $statusHandedBack

    $mdFileName = "6cd8c51c-cd40-4746-a896-72486c1d0aa6.md"
    $mdHyperlinkUrl = "https://github.com/OpenLocalizationTest/oltest/blob/60d1a9dfa8d2a640af5e7243a8e140760b642de2/e2e/$mdFileName"

    # Latest Target File (E2) - points at the source markdown file, same
    # target as the "Source File Name" hyperlink in A2.
    $e2 = $ws.Range("E2")
    $ws.Hyperlinks.Add($e2, $mdHyperlinkUrl, "", "", $mdFileName) | Out-Null
    $e2.Font.Name = "Calibri"
    $e2.Font.Size = 11
    $e2.Font.Underline = 2
    $e2.Font.Color = 15570276
    $e2.Font.Bold = $false
    $e2.Font.Italic = $false

    # Latest Handback File (F2) - points at the handed-off xlf file, same
    # target as the "Latest Handoff File" hyperlink in C2.
    $f2 = $ws.Range("F2")
    $ws.Hyperlinks.Add($f2, $XlfHyperlinkUrl, "", "", $XlfFileName) | Out-Null
    $f2.Font.Name = "Calibri"
    $f2.Font.Size = 11
    $f2.Font.Underline = 2
    $f2.Font.Color = 15570276
    $f2.Font.Bold = $false
    $f2.Font.Italic = $false

    # Latest Handback DateTime (G2)
    $ws.Range("G2").Value = $HandbackDateTime

    # Handoff Reason (H2) - the package is now included in the handback.
    $ws.Range("H2").Value = "Include"
}

Update-LanguageSheet "zh-cn" `
    "6cd8c51c-cd40-4746-a896-72486c1d0aa6.593ff4372cd824d5177f0b991a26452b69812583.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f8f7038a6aff06348199f2d50252b61fbb910ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6cd8c51c-cd40-4746-a896-72486c1d0aa6.593ff4372cd824d5177f0b991a26452b69812583.zh-cn.xlf" `
    "2016-01-18 02:23:02"

Update-LanguageSheet "de-de" `
    "6cd8c51c-cd40-4746-a896-72486c1d0aa6.593ff4372cd824d5177f0b991a26452b69812583.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/758874df52492f4410a7444d6897376aa9c46922/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6cd8c51c-cd40-4746-a896-72486c1d0aa6.593ff4372cd824d5177f0b991a26452b69812583.de-de.xlf" `
    "2016-01-18 02:23:31"

$wb.Save()
